$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-03 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-04 Sunday", 2) | Out-Null
$d.Content.Find.Execute("11÷5=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "43÷8=5, 3", 2) | Out-Null
$d.Content.Find.Execute("28÷3=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷8=5, 0", 2) | Out-Null
$d.Content.Find.Execute("68÷8=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "95÷5=19, 0", 2) | Out-Null
$d.Content.Find.Execute("21÷2=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=7, 6", 2) | Out-Null
$d.Content.Find.Execute("98÷4=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "46÷9=5, 1", 2) | Out-Null
$d.Content.Find.Execute("64÷7=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "76÷5=15, 1", 2) | Out-Null
$d.Content.Find.Execute("88÷8=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷6=13, 1", 2) | Out-Null
$d.Content.Find.Execute("11÷7=1, 4", $true, $false, $false, $false, $false, $true, 1, $false, "29÷7=4, 1", 2) | Out-Null
$d.Content.Find.Execute("82÷3=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷7=2, 0", 2) | Out-Null
$d.Content.Find.Execute("28÷8=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$d.Content.Find.Execute("93÷2=46, 1", $true, $false, $false, $false, $false, $true, 1, $false, "41÷2=20, 1", 2) | Out-Null
$d.Content.Find.Execute("22÷2=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=7, 6", 2) | Out-Null
$d.Content.Find.Execute("20÷6=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "50÷3=16, 2", 2) | Out-Null
$d.Content.Find.Execute("17÷9=1, 8", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=22, 0", 2) | Out-Null
$d.Content.Find.Execute("49÷4=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "69÷6=11, 3", 2) | Out-Null
$d.Content.Find.Execute("48÷3=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "46÷9=5, 1", 2) | Out-Null
$d.Content.Find.Execute("67÷6=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷4=19, 2", 2) | Out-Null
$d.Content.Find.Execute("61÷2=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "90÷5=18, 0", 2) | Out-Null
$d.Content.Find.Execute("19÷5=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "21÷3=7, 0", 2) | Out-Null
$d.Content.Find.Execute("11÷3=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "14÷3=4, 2", 2) | Out-Null
$d.Content.Find.Execute("98÷9=10, 8", $true, $false, $false, $false, $false, $true, 1, $false, "85÷2=42, 1", 2) | Out-Null
$d.Content.Find.Execute("60÷6=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷7=2, 5", 2) | Out-Null
$d.Content.Find.Execute("47÷6=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "78÷6=13, 0", 2) | Out-Null
$d.Content.Find.Execute("99÷5=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "18÷4=4, 2", 2) | Out-Null
$d.Content.Find.Execute("23÷7=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "76÷6=12, 4", 2) | Out-Null

Write-Host "Replacements complete"